$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("29-01-2024", "14:07:59", "registro_test2", 2.33),
    @("29-01-2024", "14:08:04", "test_form2", 5.24),
    @("29-01-2024", "14:09:25", "registro_test2", 2.34),
    @("29-01-2024", "14:09:31", "test_form2", 5.2)
)

$startRow = 8
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
